$wb = $excel.ActiveWorkbook

# --- weibull ---
$ws = $wb.Worksheets.Item("weibull")
$ws.Range("B2").Value = -2.56891849151363
$ws.Range("C2").Value = 0.0953207574839895
$ws.Range("B3").Value = -0.0434775808987697
$ws.Range("C3").Value = 0.101936547928577

# --- lognormal ---
$ws = $wb.Worksheets.Item("lognormal")
$ws.Range("B2").Value = 2.3420142874128
$ws.Range("C2").Value = 0.215810175558297
$ws.Range("B3").Value = -1.08629849598012
$ws.Range("C3").Value = 0.13211005024603

# --- llogis ---
$ws = $wb.Worksheets.Item("llogis")
$ws.Range("B2").Value = -2.08438385665272
$ws.Range("C2").Value = 0.0833283834066869
$ws.Range("B3").Value = 2.00821523370502
$ws.Range("C3").Value = 0.264693199286976

# --- gompertz ---
$ws = $wb.Worksheets.Item("gompertz")
$ws.Range("B2").Value = -2.25929198323504
$ws.Range("C2").Value = 0.08253900642319
$ws.Range("B3").Value = -0.0299328683024798
$ws.Range("C3").Value = 0.0095122801950541

# --- weibull cov ---
$ws = $wb.Worksheets.Item("weibull cov")
$ws.Range("A2").Value = 0.00908604680732155
$ws.Range("B2").Value = -0.00562530671538125
$ws.Range("A3").Value = -0.00562530671538125
$ws.Range("B3").Value = 0.0103910598035951

# --- lognormal cov ---
$ws = $wb.Worksheets.Item("lognormal cov")
$ws.Range("A2").Value = 0.0465740318745031
$ws.Range("B2").Value = -0.0269241863120591
$ws.Range("A3").Value = -0.0269241863120591
$ws.Range("B3").Value = 0.0174530653760086

# --- llogis cov ---
$ws = $wb.Worksheets.Item("llogis cov")
$ws.Range("A2").Value = 0.00694361948117181
$ws.Range("B2").Value = 0.01381622202849
$ws.Range("A3").Value = 0.01381622202849
$ws.Range("B3").Value = 0.0700624897487748

# --- gompertz cov ---
$ws = $wb.Worksheets.Item("gompertz cov")
$ws.Range("A2").Value = 0.0068126875813274
$ws.Range("B2").Value = -0.000146717072827743
$ws.Range("A3").Value = -0.000146717072827743
$ws.Range("B3").Value = 0.0000904834745092185
